$d = $word.ActiveDocument

# Target color 00B0F0 (RGB 0,176,240) expressed as the BGR integer Word uses
$blue = 15773696

# The two task paragraphs that become "done" (storage / wishlist finished)
# get their text (and paragraph mark) colored blue:
#   "*Criar função para 'favoritar' produtos - DB"
#   "*Passar classes para módulos e ver o que pode ser reaproveitado"
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if (($t -like "*favoritar*produtos*") -or ($t -like "*Passar classes para*")) {
        $p.Range.Font.Color = $blue
    }
}
